$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'275.99"
$ws.Range("E2").Value = "'-0.87%"
$ws.Range("D3").Value = "'27.33"
$ws.Range("E3").Value = "'1.63%"
$ws.Range("E4").Value = "'-2.15%"
$ws.Range("E5").Value = "'-0.71%"
$ws.Range("D6").Value = "'6.953"
$ws.Range("E6").Value = "'-0.65%"
$ws.Range("D7").Value = "'1.330"
$ws.Range("E7").Value = "'9.67%"
$ws.Range("D8").Value = "'0.8788"
$ws.Range("E8").Value = "'-0.88%"
$ws.Range("D9").Value = "'0.1521"
$ws.Range("E9").Value = "'1.99%"
$ws.Range("D10").Value = "'0.05061"
$ws.Range("E10").Value = "'-3.97%"
$ws.Range("D11").Value = "'0.07495"
$ws.Range("E11").Value = "'1.08%"
$ws.Range("D12").Value = "'0.02991"
$ws.Range("E12").Value = "'-4.78%"
$ws.Range("D13").Value = "'0.09028"
$ws.Range("E13").Value = "'-0.41%"
$ws.Range("D14").Value = "'0.001561"
$ws.Range("E14").Value = "'0.10%"
$ws.Range("D15").Value = "'0.0006407"
$ws.Range("E15").Value = "'0.98%"
$ws.Range("D16").Value = "'0.005874"
$ws.Range("E16").Value = "'-2.82%"
$ws.Range("D17").Value = "'3.452"
$ws.Range("E17").Value = "'-1.14%"
$ws.Range("D18").Value = "'3.302"
$ws.Range("E18").Value = "'-1.63%"
$ws.Range("D19").Value = "'2.284"
$ws.Range("E21").Value = "'1.61%"
$ws.Range("D22").Value = "'3.961"
$ws.Range("E22").Value = "'1.16%"
$ws.Range("D23").Value = "'0.04416"
$ws.Range("E23").Value = "'1.72%"
$ws.Range("D24").Value = "'0.001174"
$ws.Range("E24").Value = "'-0.66%"
$ws.Range("D25").Value = "'0.003867"
$ws.Range("E25").Value = "'5.16%"
$ws.Range("E26").Value = "'-0.09%"
$ws.Range("D27").Value = "'0.0001936"
$ws.Range("E27").Value = "'19.63%"
$ws.Range("D40").Value = "'0.04177"
$ws.Range("E40").Value = "'2.53%"
$ws.Range("D41").Value = "'0.006840"
$ws.Range("E41").Value = "'2.72%"
$ws.Range("E42").Value = "'0.43%"
$ws.Range("E43").Value = "'-14.06%"
$ws.Range("E44").Value = "'-9.98%"
$ws.Range("D45").Value = "'0.00005178"
$ws.Range("E45").Value = "'-1.71%"
$ws.Range("D47").Value = "'0.02299"
$ws.Range("E47").Value = "'8.26%"
